$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the previous filled row (35) onto the new row (38)
# so the new cells pick up the same fill/border/number-format styles.
$ws.Range("C35:G35").Copy()
$ws.Range("C38:G38").PasteSpecial(-4122)

# Fill in the new bill-detail entry:
# 支出 300 2018-05-10 生活费 生活费(5/11-5/20)
$ws.Range("C38").Value = "支出"
$ws.Range("D38").Value = 300
$ws.Range("E38").Value = [DateTime]"2018-05-10"
$ws.Range("F38").Value = "生活费"
$ws.Range("G38").Value = "生活费(5/11-5/20)"

# Update the view state to match (scroll position + active selection).
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I37").Select()
